$wb = $excel.ActiveWorkbook

# The original active/selected sheet (TraitDelivery_TraitID, 4th tab).
$ws4 = $wb.Worksheets.Item(4)

# Duplicate it to get an identical copy of all formatting/column widths,
# then turn that copy into the new "TraitCon_TraitID" sheet (trait-level
# conversion lookup), placed right after it (i.e. at the end).
$ws4.Copy($null, $ws4) | Out-Null
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "TraitCon_TraitID"

# The template sheet has 5 data rows; the new sheet only needs 4 (header + 3).
$newSheet.Rows.Item(5).Delete() | Out-Null

# Fill in the trait-level conversion mapping rows. Columns A and B are
# written together, row-by-row for rows 3 and 4 first, and the A2 cell
# last, so new shared strings are created in the same relative order as
# the source workbook.
$newSheet.Range("A3").Value = "total_click_based_conversion"
$newSheet.Range("A4").Value = "total_view_based_conversion"
$newSheet.Range("B3").Value = "Trait Click Based Conversions"
$newSheet.Range("B4").Value = "Trait View Based Conversions"
$newSheet.Range("A2").Value = "category_id"

$newSheet.Range("B2").Value = "Segment ID"

$newSheet.Range("C2").Value = $true
$newSheet.Range("C3").Value = $false
$newSheet.Range("C4").Value = $false

$newSheet.Range("D2").Value = "INT"
$newSheet.Range("D3").Value = "DOUBLE"
$newSheet.Range("D4").Value = "DOUBLE"

$newSheet.Range("E2").Value = "MATCH"
$newSheet.Range("E3").Value = "SUBSTRING"
$newSheet.Range("E4").Value = "MATCH"

# New sheet becomes the active tab; A2 is the selected cell on it.
$newSheet.Range("A2").Select() | Out-Null
$newSheet.PageSetup.Orientation = 1

# The previously-active sheet is no longer the selected tab; its
# selection reverts to the whole-sheet default.
$ws4.Range("A1:XFD1048576").Select() | Out-Null
$ws4.PageSetup.Orientation = 1

$newSheet.Activate() | Out-Null

Write-Host "Done"
